$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    ,@(16, "45553634", "DINA LUZ CALVO RAMIREZ", "1905", 33125, 828116)
    ,@(17, "40740818", "LEIDY PAOLA CHILATRA MONTOYA", "1905", 33125, 828116)
    ,@(18, "1032358550", "RICARDO JOSE DIAZ VERGARA", "1905", 33125, 828116)
    ,@(19, "1128052485", "JULIO CESAR LUNA CASANOVA", "1905", 33125, 737717)
    ,@(20, "1051449332", "DAHISANIS MUENTES HERNANDEZ", "1905", 33125, 828116)
    ,@(21, "45553634", "DINA LUZ CALVO RAMIREZ", "1906", 33125, 828116)
    ,@(22, "40740818", "LEIDY PAOLA CHILATRA MONTOYA", "1906", 33125, 828116)
    ,@(23, "1032358550", "RICARDO JOSE DIAZ VERGARA", "1906", 33125, 828116)
    ,@(24, "1128052485", "JULIO CESAR LUNA CASANOVA", "1906", 33125, 737717)
    ,@(25, "1051449332", "DAHISANIS MUENTES HERNANDEZ", "1906", 33125, 828116)
    ,@(26, "45553634", "DINA LUZ CALVO RAMIREZ", "1907", 33125, 828116)
    ,@(27, "40740818", "LEIDY PAOLA CHILATRA MONTOYA", "1907", 33125, 828116)
    ,@(28, "1032358550", "RICARDO JOSE DIAZ VERGARA", "1907", 33125, 828116)
    ,@(29, "1128052485", "JULIO CESAR LUNA CASANOVA", "1907", 33125, 737717)
    ,@(30, "1051449332", "DAHISANIS MUENTES HERNANDEZ", "1907", 33125, 828116)
    ,@(31, "45553634", "DINA LUZ CALVO RAMIREZ", "1908", 33125, 828116)
    ,@(32, "40740818", "LEIDY PAOLA CHILATRA MONTOYA", "1908", 33125, 828116)
    ,@(33, "1032358550", "RICARDO JOSE DIAZ VERGARA", "1908", 33125, 828116)
    ,@(34, "1128052485", "JULIO CESAR LUNA CASANOVA", "1908", 33125, 737717)
    ,@(35, "1051449332", "DAHISANIS MUENTES HERNANDEZ", "1908", 33125, 828116)
    ,@(36, "45553634", "DINA LUZ CALVO RAMIREZ", "1909", 33125, 828116)
    ,@(37, "40740818", "LEIDY PAOLA CHILATRA MONTOYA", "1909", 33125, 828116)
    ,@(38, "1032358550", "RICARDO JOSE DIAZ VERGARA", "1909", 33125, 828116)
    ,@(39, "1128052485", "JULIO CESAR LUNA CASANOVA", "1909", 33125, 737717)
    ,@(40, "1051449332", "DAHISANIS MUENTES HERNANDEZ", "1909", 33125, 828116)
    ,@(41, "45553634", "DINA LUZ CALVO RAMIREZ", "1910", 33125, 828116)
    ,@(42, "40740818", "LEIDY PAOLA CHILATRA MONTOYA", "1910", 33125, 828116)
    ,@(43, "1032358550", "RICARDO JOSE DIAZ VERGARA", "1910", 33125, 828116)
    ,@(44, "1128052485", "JULIO CESAR LUNA CASANOVA", "1910", 33125, 737717)
    ,@(45, "1051449332", "DAHISANIS MUENTES HERNANDEZ", "1910", 33125, 828116)
    ,@(46, "45553634", "DINA LUZ CALVO RAMIREZ", "1911", 26500, 828116)
    ,@(47, "40740818", "LEIDY PAOLA CHILATRA MONTOYA", "1911", 26500, 828116)
    ,@(48, "1032358550", "RICARDO JOSE DIAZ VERGARA", "1911", 26500, 828116)
    ,@(49, "1128052485", "JULIO CESAR LUNA CASANOVA", "1911", 23607, 737717)
    ,@(50, "1051449332", "DAHISANIS MUENTES HERNANDEZ", "1911", 26500, 828116)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}
